$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5
$ws.Range("A5").Value = "f34630c2-31aa-11f0-8dd5-fa163ee583d0"
$ws.Range("B5").Value = "15.05.2025"
$ws.Range("E5").Value = "loris.marino@students.fhnw.ch"

# F5 holds a phone number that looks numeric; force text storage (like the
# other rows' phone cells) without leaving a residual non-default style.
$ws.Range("F5").NumberFormat = "@"
$ws.Range("F5").Value = "41787153999"
$ws.Range("F5").Style = "Normal"

$ws.Range("G5").Value = "Hans"
$ws.Range("H5").Value = "Volter"
$ws.Range("I5").Value = "Many have tried to interfere with my work. You think you are somehow better?"

# Row 6
$ws.Range("A6").Value = "140a4b37-31ac-11f0-8dd5-fa163ee583d0"
$ws.Range("B6").Value = "15.05.2025"
$ws.Range("E6").Value = "loris.marino@students.fhnw.ch"

$ws.Range("F6").NumberFormat = "@"
$ws.Range("F6").Value = "41787153999"
$ws.Range("F6").Style = "Normal"

$ws.Range("G6").Value = "Loris"
$ws.Range("H6").Value = "Mariño"
$ws.Range("I6").Value = "Feedback Feedback DASSSSSSSSSSSSSSSSSSSSSSSSSSSSSSSSSasda sdjasdi asdjasd asd adihsdifhse IFESUFH SEIUFHESUFH SEIFSEIFH SEIIUFHUD SDLJHSDK JVHKJ YDkjefsdkjfh sifhsduf hsfhefsdv90shgw89ej f9sdv9'ah erduvi sdfgnraivndf ipdufn ipjdfn era9u a9erngijsdvnpiagh 9'a4eugaidrug iardug 'a8eugiuergn rghaergnuaer giuaeug aegaerug hae9rghae9r greg oyidjk<setm nbglkdifjvpiuvpiruagh98a hyidnvpivjpiaegh a iaug aerg raiugidugh a'ivfdipghea"
